$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.607.67'
$ws.Range("E2").Value = '  -3.69%  '

$ws.Range("D3").Value = '2.513.56'
$ws.Range("E3").Value = '  -4.95%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '579.01'
$ws.Range("E5").Value = '  -2.04%  '

$ws.Range("D6").Value = '167.36'
$ws.Range("E6").Value = '  -4.37%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  -0.63%  '

$ws.Range("D9").Value = '2.513.78'
$ws.Range("E9").Value = '  -4.91%  '

$ws.Range("E10").Value = '  -6.39%  '

$ws.Range("E11").Value = '  -0.51%  '

$ws.Range("D12").Value = '0.342'
$ws.Range("E12").Value = '  -4.01%  '

$ws.Range("D13").Value = '4.87'
$ws.Range("E13").Value = '  -1.97%  '

$ws.Range("D14").Value = '2.964.57'
$ws.Range("E14").Value = '  -5.29%  '

$ws.Range("D15").Value = '69.469.04'
$ws.Range("E15").Value = '  -3.72%  '

$ws.Range("E16").Value = '  -5.33%  '

$ws.Range("D17").Value = '24.96'
$ws.Range("E17").Value = '  -4.04%  '

$ws.Range("D18").Value = '2.510.26'
$ws.Range("E18").Value = '  -3.83%  '

$ws.Range("E19").Value = '  -6.78%  '

$ws.Range("D20").Value = '7.79'
$ws.Range("E20").Value = '  -2.59%  '

$ws.Range("D21").Value = '351.72'
$ws.Range("E21").Value = '  -4.90%  '

$ws.Range("E22").Value = '  -4.54%  '

$ws.Range("E23").Value = '  -3.30%  '

$ws.Range("E24").Value = '  +0.05%  '

$ws.Range("D25").Value = '69.29'
$ws.Range("E25").Value = '  -2.95%  '

$ws.Range("E26").Value = '  -5.35%  '

$ws.Range("D27").Value = '9.07'
$ws.Range("E27").Value = '  -6.26%  '

$ws.Range("D28").Value = '2.642.31'
$ws.Range("E28").Value = '  -5.01%  '

$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.27%  '

$ws.Range("D30").Value = '0.0₃0908'
$ws.Range("E30").Value = '  -4.75%  '

$ws.Range("D31").Value = '7.90'
$ws.Range("E31").Value = '  -1.85%  '

$ws.Range("D32").Value = '479.61'
$ws.Range("E32").Value = '  -3.67%  '

$ws.Range("E33").Value = '  +0.83%  '

$ws.Range("E34").Value = '  -2.77%  '

$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("E36").Value = '  -0.89%  '

$ws.Range("D37").Value = '152.60'

$ws.Range("D38").Value = '18.90'
$ws.Range("E38").Value = '  -0.04%  '

$ws.Range("D39").Value = '18.61'
$ws.Range("E39").Value = '  -3.88%  '

$ws.Range("E40").Value = '  -0.04%  '

$ws.Range("D41").Value = '4.79'
$ws.Range("E41").Value = '  -2.57%  '

$ws.Range("E42").Value = '  -2.80%  '

$ws.Range("D43").Value = '1.62'
$ws.Range("E43").Value = '  -6.17%  '

$ws.Range("E44").Value = '  -13.37%  '

$ws.Range("E45").Value = '  -8.23%  '

$ws.Range("D46").Value = '38.18'
$ws.Range("E46").Value = '  -2.52%  '

$ws.Range("D47").Value = '144.45'
$ws.Range("E47").Value = '  -6.16%  '

$ws.Range("D48").Value = '3.56'
$ws.Range("E48").Value = '  -3.16%  '

$ws.Range("D49").Value = '0.533'
$ws.Range("E49").Value = '  -3.32%  '

$ws.Range("E50").Value = '  -4.72%  '

$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '0.588'
$ws.Range("E51").Value = '  -1.96%  '
